$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values such as "301.02" or
# "43.039.17" are stored as literal text (matching the source inlineStr cells)
# instead of being auto-converted to numbers by Excel's type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '43.039.17'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '2.304.04'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '301.02'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = '98.35'
$ws.Range('E6').Value = '  -1.39%  '
$ws.Range('E7').Value = '  +2.06%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('D10').Value = '36.31'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '17.89'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').Value = '6.82'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = '2.662.68'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('D16').Value = '2.310.06'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').Value = '42.996.86'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '12.62'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = '0.0₃0910'
$ws.Range('D21').Value = '6.12'
$ws.Range('E21').Value = '  -1.54%  '
$ws.Range('D22').Value = '68.35'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '242.38'
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').Value = '25.25'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').Value = '166.92'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D31').Value = '9.10'
$ws.Range('D32').Value = '33.23'
$ws.Range('E32').Value = '  -3.28%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('D36').Value = '17.77'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').Value = '2.002.48'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('D46').Value = '10.21'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').Value = '17.54'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('D49').Value = '53.72'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('D50').Value = '2.528.40'
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').Value = '72.79'
$ws.Range('E51').Value = '  -4.45%  '
